{"js": "// Remove the \"Do poe mapping system...\" / \"In the hideout...\" /\n// \"Ideally, the hideout should be Shaper's...\" paragraphs (plus the two\n// blank paragraphs interleaved with them) that used to sit right after the\n// \"before the release of Godot4.1...\" paragraph in the \"error handling\"\n// section of the report.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph that must stay untouched.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"before the release of Godot4.1\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the Godot4.1 anchor paragraph.\");\n}\n\n// Locate the last paragraph that must be removed (the \"Ideally, the\n// hideout...\" paragraph) so we know exactly where the deletion block ends.\nlet lastRemoveIndex = -1;\nfor (let i = anchorIndex + 1; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ideally, the hideout should be Shaper\") !== -1) {\n    lastRemoveIndex = i;\n    break;\n  }\n}\n\nif (lastRemoveIndex === -1) {\n  throw new Error(\"Could not locate the trailing paragraph to remove.\");\n}\n\n// Delete every paragraph strictly between the anchor and (inclusive of) the\n// \"Ideally, the hideout...\" paragraph.\nfor (let i = anchorIndex + 1; i <= lastRemoveIndex; i++) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Do poe mapping system...\" / \"In the hideout...\" /\n# \"Ideally, the hideout should be Shaper's...\" paragraphs (plus the two\n# blank paragraphs interleaved with them) that used to sit right after the\n# \"before the release of Godot4.1...\" paragraph in the report.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph that must stay untouched.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*before the release of Godot4.1*\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the Godot4.1 anchor paragraph.\"\n}\n\n# Locate the last paragraph that must be removed (the \"Ideally, the\n# hideout...\" paragraph) so we know exactly where the deletion block ends.\n$lastRemoveIndex = -1\nfor ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Ideally, the hideout should be Shaper*\") {\n        $lastRemoveIndex = $i\n        break\n    }\n}\nif ($lastRemoveIndex -eq -1) {\n    throw \"Could not locate the trailing paragraph to remove.\"\n}\n\n# Delete the whole block (everything strictly after the anchor through, and\n# including, the \"Ideally, the hideout...\" paragraph) in one go.\n$startPara = $d.Paragraphs.Item($anchorIndex + 1)\n$endPara = $d.Paragraphs.Item($lastRemoveIndex)\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n"}
